$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.702.45'
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").Value = '2.173.52'
$ws.Range("E3").Value = '  -2.70%  '
$ws.Range("E4").Value = '  -0.11%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.24'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -1.95%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.605'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -3.49%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.71'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -1.78%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  -3.07%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.06'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -5.25%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -4.90%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.55'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  -3.56%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0999'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -2.89%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.68'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("D15").Value = '2.502.16'
$ws.Range("E15").Value = '  -2.60%  '
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.36'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").Value = '2.155.57'
$ws.Range("E17").Value = '  -2.78%  '
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.781'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -6.73%  '
$ws.Range("D19").Value = '41.625.80'
$ws.Range("E19").Value = '  -1.13%  '
$ws.Range("E20").Value = '  -2.22%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.18'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -3.57%  '
$ws.Range("E22").Value = '  -6.64%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.09'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -11.00%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '226.34'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("E25").Value = '  +0.16%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.75'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -5.72%  '
$ws.Range("E28").Value = '  -9.94%  '
$ws.Range("E29").Value = '  -3.56%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.88'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +2.24%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.08'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  -5.17%  '
$ws.Range("E32").Value = '  -3.50%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.76'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +9.97%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0773'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  -3.89%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.29'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -6.01%  '
$ws.Range("E36").Value = '  -4.06%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.30'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -1.60%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.102'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -6.95%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0307'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +0.90%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.09'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -8.53%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.09'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -1.82%  '
$ws.Range("E42").Value = '  -5.80%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.31'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -8.08%  '
$ws.Range("E44").Value = '  -4.96%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.46'
$ws.Range("D45").Style = $origStyle
$ws.Range("E46").Value = '  -3.51%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.32'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -6.82%  '
$ws.Range("E48").Value = '  -4.97%  '
$ws.Range("E49").Value = '  -5.31%  '
$ws.Range("E50").Value = '  -6.05%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.63'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -2.02%  '
